# Revert "adding term 2.0.0"
# - Delete worksheet "Include from FSIII 2"
# - Restore Version, Date, Contact values on Metadata sheet
# - Restore value in Include from FSIII sheet C2 cell

$wb = $excel.ActiveWorkbook

# Update Metadata sheet values
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# Update Include from FSIII sheet value
$inc = $wb.Worksheets.Item("Include from FSIII")
$inc.Range("C2").Value = "E"

# Delete the extra "Include from FSIII 2" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Include from FSIII 2").Delete()
$excel.DisplayAlerts = $true
